$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text so values like "1.004" are not
# auto-converted into numbers by Excel type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.043.65"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.869.95"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "312.84"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "0.5060"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("D8").Value = "0.3816"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").Value = "0.08311"
$ws.Range("E9").Value = "  -9.78%  "

$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -1.43%  "

$ws.Range("D11").Value = "41.43"
$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("D12").Value = "6.208"
$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("D13").Value = "1.869.92"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").Value = "20.49"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "7.191"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").Value = "90.75"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").Value = "0.06634"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "17.93"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D22").Value = "6.027"
$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").Value = "28.078.16"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("D26").Value = "2.577"
$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").Value = "2.087.23"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").Value = "156.72"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").Value = "20.61"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("D30").Value = "125.59"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("D31").Value = "0.1054"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").Value = "1.044"
$ws.Range("E32").Value = "  -3.21%  "

$ws.Range("D33").Value = "5.600"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "3.602"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "9.713"
$ws.Range("E35").Value = "  +2.73%  "

$ws.Range("D36").Value = "0.02450"
$ws.Range("E36").Value = "  +2.06%  "

$ws.Range("D37").Value = "0.06564"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").Value = "0.2164"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("E39").Value = "  -0.73%  "

$ws.Range("D40").Value = "0.6455"
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("E41").Value = "  -7.51%  "

$ws.Range("D42").Value = "11.30"
$ws.Range("E42").Value = "  -2.26%  "

$ws.Range("D43").Value = "4.878"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("D44").Value = "0.6140"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("D45").Value = "13.08"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").Value = "1.293"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "2.009"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").Value = "1.215"
$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("D50").Value = "121.41"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").Value = "80.19"
$ws.Range("E51").Value = "  +1.62%  "

# Restore the default (Normal) style on column D so no stray number format
# metadata is left attached to the cells (matches original workbook styling).
$ws.Range("D2:D51").Style = "Normal"
